# Applies the "Fixed General values to test diferent solution" edit:
# toggles a set of "x" marker cells across the Año1/Año2/Año3 sheets of the
# horariosProhibidos workbook, then restores the per-sheet cell selections
# and leaves Año3 as the active tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Año1")
$ws2 = $wb.Worksheets.Item("Año2")
$ws3 = $wb.Worksheets.Item("Año3")

# ---------------------------------------------------------------------
# Año1 data changes
# ---------------------------------------------------------------------
$ws1.Range("H5").Value = "x"

$ws1.Range("E6").Value = "x"
$ws1.Range("G6").Value = "x"
$ws1.Range("H6").Value = "x"

$ws1.Range("E7").Value = "x"
$ws1.Range("G7").Value = "x"
$ws1.Range("H7").Value = "x"

$ws1.Range("H8").Value = "x"

$ws1.Range("H10").Value = "x"

$ws1.Range("F11").Value = "x"
$ws1.Range("H11").Value = "x"

$ws1.Range("F12").Value = "x"
$ws1.Range("H12").Value = "x"

$ws1.Range("H13").Value = "x"

$ws1.Range("D22:H22").ClearContents() | Out-Null

# ---------------------------------------------------------------------
# Año2 data changes
# ---------------------------------------------------------------------
$ws2.Range("D5").Value = "x"
$ws2.Range("F5").ClearContents() | Out-Null
$ws2.Range("H5").Value = "x"

$ws2.Range("D6").Value = "x"
$ws2.Range("F6").ClearContents() | Out-Null
$ws2.Range("H6").Value = "x"

$ws2.Range("F7").ClearContents() | Out-Null
$ws2.Range("H7").Value = "x"

$ws2.Range("F9").Value = "x"

$ws2.Range("F10").ClearContents() | Out-Null

$ws2.Range("F11").ClearContents() | Out-Null

$ws2.Range("D12").Value = "x"
$ws2.Range("H12").Value = "x"

$ws2.Range("D13").Value = "x"
$ws2.Range("H13").Value = "x"

# ---------------------------------------------------------------------
# Año3 data changes
# ---------------------------------------------------------------------
$ws3.Range("G5").ClearContents() | Out-Null
$ws3.Range("G6").ClearContents() | Out-Null
$ws3.Range("G7").ClearContents() | Out-Null
$ws3.Range("G8").ClearContents() | Out-Null

$ws3.Range("G10").ClearContents() | Out-Null
$ws3.Range("G11").ClearContents() | Out-Null

$ws3.Range("D12").Value = "x"
$ws3.Range("G12").ClearContents() | Out-Null
$ws3.Range("H12").Value = "x"

$ws3.Range("D13").Value = "x"
$ws3.Range("G13").ClearContents() | Out-Null
$ws3.Range("H13").Value = "x"

$ws3.Range("D14").Value = "x"
$ws3.Range("H14").Value = "x"

$ws3.Range("F15").Value = "x"

$ws3.Range("F16").Value = "x"

$ws3.Range("E19").Value = "x"
$ws3.Range("G19").Value = "x"

$ws3.Range("E20").Value = "x"
$ws3.Range("G20").Value = "x"

$ws3.Range("E21").Value = "x"
$ws3.Range("G21").Value = "x"

# ---------------------------------------------------------------------
# Restore per-sheet selections. Selecting a range switches the active
# sheet/tab, so Año3 (the desired final active tab) is selected last.
# ---------------------------------------------------------------------
$ws1.Range("G18:G19").Select() | Out-Null
$ws2.Range("H21").Select() | Out-Null
$ws3.Range("H14").Select() | Out-Null
